$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.955.51"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "2.883.31"
$ws.Range("E3").Value = "  -1.53%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'587.28"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "'138.48"
$ws.Range("E6").Value = "  -5.92%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -3.41%  "
$ws.Range("D9").Value = "'6.88"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("E10").Value = "  -5.12%  "
$ws.Range("E11").Value = "  -3.44%  "
$ws.Range("D12").Value = "'0.0000217"
$ws.Range("E12").Value = "  -4.40%  "
$ws.Range("D13").Value = "'32.17"
$ws.Range("E13").Value = "  -4.57%  "
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "3.358.99"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "60.878.47"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "2.884.76"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("E18").Value = "  -3.77%  "
$ws.Range("D19").Value = "'423.42"
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").Value = "'13.23"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("D21").Value = "'0.651"
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("D23").Value = "'79.76"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("D24").Value = "'10.36"
$ws.Range("E24").Value = "  -4.89%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "'2.05"
$ws.Range("E26").Value = "  -7.48%  "
$ws.Range("E27").Value = "  -5.18%  "
$ws.Range("E28").Value = "  -3.81%  "
$ws.Range("D29").Value = "'2.06"
$ws.Range("E29").Value = "  -9.31%  "
$ws.Range("E30").Value = "  -6.33%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'25.55"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("E33").Value = "  -5.90%  "
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("D35").Value = "'0.967"
$ws.Range("E35").Value = "  -4.84%  "
$ws.Range("E36").Value = "  -4.42%  "
$ws.Range("D37").Value = "'49.00"
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("D38").Value = "'2.78"
$ws.Range("E38").Value = "  -7.94%  "
$ws.Range("E39").Value = "  -4.79%  "
$ws.Range("D40").Value = "'8.32"
$ws.Range("E40").Value = "  -3.01%  "
$ws.Range("E41").Value = "  -6.44%  "
$ws.Range("E42").Value = "  -7.03%  "
$ws.Range("D43").Value = "2.662.45"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("D44").Value = "'37.63"
$ws.Range("E44").Value = "  -7.79%  "
$ws.Range("D45").Value = "'131.71"
$ws.Range("E45").Value = "  -1.76%  "
$ws.Range("E46").Value = "  -4.69%  "
$ws.Range("D47").Value = "'347.04"
$ws.Range("E47").Value = "  -8.62%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  -4.53%  "
$ws.Range("E50").Value = "  -7.25%  "
$ws.Range("D51").Value = "'1.92"
$ws.Range("E51").Value = "  -4.90%  "
